$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: sd / Statement-non-opinion -> b / Acknowledge (Backchannel)
$ws.Range("I5").Value = "b"
$ws.Range("J5").Value = "Acknowledge (Backchannel)"

# Row 7: qy / Yes-No-Question -> sv / Statement-opinion
$ws.Range("I7").Value = "sv"
$ws.Range("J7").Value = "Statement-opinion"

# Row 17: sd / Statement-non-opinion -> % / Uninterpretable
$ws.Range("I17").Value = "%"
$ws.Range("J17").Value = "Uninterpretable"

# Row 21: sd / Statement-non-opinion -> aa / Agree/Accept
$ws.Range("I21").Value = "aa"
$ws.Range("J21").Value = "Agree/Accept"

# Row 26: b / Acknowledge (Backchannel) -> sd / Statement-non-opinion
$ws.Range("I26").Value = "sd"
$ws.Range("J26").Value = "Statement-non-opinion"
